# Applies:
#  1. Rename sheet "Cetacean Genomes Project" -> "Cetaceans Genomes Project"
#  2. Update column A (Affiliated Project) values on that sheet from
#     "Cetacean Genomes Project (CGP)" -> "Cetaceans Genomes Project (CGP)"
#  3. Fix the Bogotá address row on the "Wise Ancestors" sheet and
#     populate its previously-blank Latitude/Longitude cells.

$wb = $excel.ActiveWorkbook

# --- 1 & 2: Cetacean(s) Genomes Project sheet -------------------------------
$ws = $wb.Worksheets.Item("Cetacean Genomes Project")
$ws.Name = "Cetaceans Genomes Project"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value() -eq "Cetacean Genomes Project (CGP)") {
        $cell.Value = "Cetaceans Genomes Project (CGP)"
    }
}

# --- 3: Wise Ancestors sheet, row 3 (Bogotá address) ------------------------
$wa = $wb.Worksheets.Item("Wise Ancestors")
$wa.Range("E3").Value = "Calle 72 - 65 Piso 7, Chapinero, Bogotá, Cundinamarca, Colombia"
$wa.Range("F3").Value = 4.6535353
$wa.Range("G3").Value = -74.05484229999999
